$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (B22) down into the new B23/B24 cells
# so the new "Masuk" cells pick up the same currency number style (s="2").
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B23:B24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 23: new transaction entry (17 Mei 2023 update)
$ws.Range("A23").Value = 45057
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 200000
$ws.Range("E23").Value = "uang santunan dukacita untuk keluargag kak safe"
$ws.Range("F23").Value = "yofandi"

# Row 24: new transaction entry
$ws.Range("A24").Value = 45063
$ws.Range("B24").Value = 70000
$ws.Range("C24").Value = 0
$ws.Range("E24").Value = "uang persembahan - reguler"
$ws.Range("F24").Value = "yofandi"

# Update the view state to match the author's saved selection/scroll position
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F24").Select() | Out-Null
